$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.279.53'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').Value = '2.059.43'
$ws.Range('E3').Value = '  +3.20%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '234.27'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = '0.612'
$ws.Range('E6').Value = '  +2.13%  '
$ws.Range('D7').Value = '58.12'
$ws.Range('E7').Value = '  +5.85%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +2.35%  '
$ws.Range('D10').Value = '58.73'
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('D11').Value = '0.0761'
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('E12').Value = '  +2.70%  '
$ws.Range('D13').Value = '2.364.65'
$ws.Range('E13').Value = '  +3.35%  '
$ws.Range('D14').Value = '14.56'
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('D15').Value = '21.09'
$ws.Range('E15').Value = '  +3.07%  '
$ws.Range('D16').Value = '0.773'
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').Value = '2.062.09'
$ws.Range('E18').Value = '  +3.36%  '
$ws.Range('D19').Value = '37.527.20'
$ws.Range('E19').Value = '  +2.72%  '
$ws.Range('E20').Value = '  +16.45%  '
$ws.Range('D21').Value = '69.53'
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('D22').Value = '0.0₃0813'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('D23').Value = '226.87'
$ws.Range('E23').Value = '  +2.22%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  +1.39%  '
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('D27').Value = '165.27'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').Value = '1.49'
$ws.Range('E28').Value = '  +11.14%  '
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('D30').Value = '19.17'
$ws.Range('E30').Value = '  +1.61%  '
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('D32').Value = '0.118'
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('E33').Value = '  +3.26%  '
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('E35').Value = '  +8.78%  '
$ws.Range('D36').Value = '4.54'
$ws.Range('E36').Value = '  +6.23%  '
$ws.Range('D37').Value = '3.37'
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').Value = '1.79'
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('E40').Value = '  +4.06%  '
$ws.Range('D41').Value = '0.0981'
$ws.Range('E41').Value = '  +3.53%  '
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('D43').Value = '4.34'
$ws.Range('E43').Value = '  +22.16%  '
$ws.Range('D44').Value = '1.454.92'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '95.42'
$ws.Range('E45').Value = '  +6.93%  '
$ws.Range('E46').Value = '  +3.85%  '
$ws.Range('E47').Value = '  +3.97%  '
$ws.Range('D48').Value = '15.80'
$ws.Range('E48').Value = '  +3.58%  '
$ws.Range('D49').Value = '1.02'
$ws.Range('E49').Value = '  +2.95%  '
$ws.Range('D50').Value = '7.24'
$ws.Range('E50').Value = '  +5.45%  '
$ws.Range('E51').Value = '  +1.91%  '
